$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-08"

# Update the shared header string cell (I1) that reads "2022 (through 09-07)"
$ws.Range("I1").Value = "2022 (through 09-08)"

# Update September row (row 10) value in column I
$ws.Range("I10").Value = 38

# Update Total row (row 14) value in column I
$ws.Range("I14").Value = 1175
